$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 5 ("Bec") with new ROI values
$ws.Range("B5").Value = 1172
$ws.Range("C5").Value = 1350
$ws.Range("D5").Value = 1667
$ws.Range("E5").Value = 1881
$ws.Range("H5").Value = 349.10000000000002

# Add new row 9 ("NiLattice") with ROI values
$ws.Range("A9").Value = "NiLattice"
$ws.Range("B9").Value = 1107
$ws.Range("C9").Value = 1417
$ws.Range("D9").Value = 1740
$ws.Range("E9").Value = 1814
$ws.Range("F9").Value = 2160
$ws.Range("G9").Value = 2560
$ws.Range("H9").Value = 349.10000000000002
